$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the redundant duplicated "mainField" column (column D) from the
# pack_rule1() decision table (rows 13-16). The data in column D was an
# exact duplicate of column C and is no longer needed now that columns
# with true conditions are matched directly with the output object.
$ws.Range("D13:D16").ClearContents()

# Restore the active selection/view state used after the edit.
$ws.Range("H19").Select()
